# Refresh the cryptos price/volume table (Coin, Link, Price, Volume(1h))
# with the latest scrape. Column A (rank index) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.761.14"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").Value = "3.053.65"
$ws.Range("E3").Value = "  +2.35%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'524.24"
$ws.Range("E5").Value = "  +5.47%  "

$ws.Range("D6").Value = "'141.88"
$ws.Range("E6").Value = "  +5.26%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +4.77%  "

$ws.Range("D9").Value = "'7.63"
$ws.Range("E9").Value = "  +4.96%  "

$ws.Range("E10").Value = "  +7.91%  "

$ws.Range("D11").Value = "'0.370"
$ws.Range("E11").Value = "  +4.96%  "

$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").Value = "3.573.24"

$ws.Range("D14").Value = "'26.91"
$ws.Range("E14").Value = "  +8.14%  "

$ws.Range("D15").Value = "'0.0000170"
$ws.Range("E15").Value = "  +16.56%  "

$ws.Range("D16").Value = "57.756.93"
$ws.Range("E16").Value = "  +2.18%  "

$ws.Range("D17").Value = "'6.24"
$ws.Range("E17").Value = "  +7.09%  "

$ws.Range("D18").Value = "3.059.29"
$ws.Range("E18").Value = "  +2.59%  "

$ws.Range("D19").Value = "'13.04"
$ws.Range("E19").Value = "  +5.52%  "

$ws.Range("D20").Value = "'8.17"
$ws.Range("E20").Value = "  +5.45%  "

$ws.Range("D21").Value = "'339.41"
$ws.Range("E21").Value = "  +4.49%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'0.501"
$ws.Range("E23").Value = "  +7.55%  "

$ws.Range("D24").Value = "'64.92"
$ws.Range("E24").Value = "  +5.97%  "

$ws.Range("E25").Value = "  +6.34%  "

$ws.Range("D26").Value = "0.0₃0973"
$ws.Range("E26").Value = "  +8.96%  "

$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("D28").Value = "'6.92"
$ws.Range("E28").Value = "  +6.81%  "

$ws.Range("E29").Value = "  +10.43%  "

$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "  +7.10%  "

$ws.Range("E31").Value = "  +5.34%  "

$ws.Range("D32").Value = "'21.09"
$ws.Range("E32").Value = "  +4.66%  "

# Row 33 now lists Monero
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'156.50"
$ws.Range("E33").Value = "  +0.64%  "

# Row 34 now lists NEARProtocol
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  +6.15%  "

$ws.Range("D35").Value = "'5.96"
$ws.Range("E35").Value = "  +6.41%  "

$ws.Range("D36").Value = "'1.32"
$ws.Range("E36").Value = "  +3.23%  "

$ws.Range("D37").Value = "'26.01"
$ws.Range("E37").Value = "  +11.70%  "

$ws.Range("D38").Value = "'0.0705"
$ws.Range("E38").Value = "  +3.91%  "

$ws.Range("D39").Value = "3.088.77"
$ws.Range("E39").Value = "  +2.45%  "

$ws.Range("D40").Value = "'37.71"
$ws.Range("E40").Value = "  +3.34%  "

$ws.Range("E41").Value = "  +9.00%  "

# Row 42 now lists FirstDigitalUSD
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43 now lists Stacks
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.48"
$ws.Range("E43").Value = "  +5.47%  "

# Row 44 now lists Mantle
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.664"
$ws.Range("E44").Value = "  +4.25%  "

# Row 45 now lists Maker
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.330.30"
$ws.Range("E45").Value = "  +4.71%  "

$ws.Range("E46").Value = "  +4.21%  "

$ws.Range("E47").Value = "  +4.46%  "

$ws.Range("E48").Value = "  +4.52%  "

$ws.Range("D49").Value = "'6.05"
$ws.Range("E49").Value = "  +4.59%  "

$ws.Range("D50").Value = "'20.20"
$ws.Range("E50").Value = "  +5.87%  "

$ws.Range("E51").Value = "  +6.10%  "
